$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9025498628616333
$ws.Range("B1").Value = 1.667898416519165
$ws.Range("C1").Value = 4.300154209136963
$ws.Range("D1").Value = 2.848291158676147
$ws.Range("E1").Value = 0.6086277961730957
